# Export des ontologies, entités de dbpedia de Mtab.
# - Remove column D (birthPlace header + its two data values)
# - Append a new row 4 (Y._D._Tiwari / Shillong Meghalaya / Father)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column D entirely (header + D2/D3 data) - also re-tightens the
# sheet dimension from A1:D3 down to A1:C3.
$ws.Range("D1:D3").Delete()

# Add the new 4th row of data (A1:C4 is the final used range).
$ws.Range("A4").Value = "http://dbpedia.org/resource/Y._D._Tiwari"
$ws.Range("B4").Value = "http://dbpedia.org/resource/Shillong http://dbpedia.org/resource/Meghalaya "
$ws.Range("C4").Value = "http://dbpedia.org/resource/Father"
